# data cleanup continued in player_per_game_df
# Insert a new row for "LeBron James" into the alphabetically-sorted
# player/award pivot table on Sheet1, between "Larry Bird" (row 22)
# and "Magic Johnson" (existing row 23, which shifts down to row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 23:35 down to 24:36, leaving a blank row 23 for the new entry.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row.
$ws.Cells.Item(23, 1).Value2 = "LeBron James"
$ws.Cells.Item(23, 2).Value2 = 10
